$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new period column "Agosto.2021" in column BH, one column to the
# right of the existing last column BG ("Mayo.2021"). The new column's
# values repeat the latest known figures from column BG (same as BF->BG
# pattern already present in the sheet), for all 18 data rows (2-19).

# Copy the header formatting (bold font, border, centered alignment) from
# the previous header cell (BG1) onto the new header cell (BH1), then set
# its text.
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BH1").Value = "Agosto.2021"

for ($r = 2; $r -le 19; $r++) {
    $srcCell = $ws.Cells.Item($r, 59)   # column BG
    $dstCell = $ws.Cells.Item($r, 60)   # column BH
    $dstCell.Value = $srcCell.Value2
}
